$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 headers: shift values C1->E1, D1->C1, E1->D1 (new order: prediction, rejection-f, max)
$ws.Range("C1").Value = "prediction"
$ws.Range("D1").Value = "rejection-f"
$ws.Range("E1").Value = "max"

# Row 2 values
$ws.Range("C2").Value = "s__CAG-964 sp000435335"
$ws.Range("D2").Value = "s__CAG-964 sp000435335"
$ws.Range("E2").Value = 0.9997373809691441
